# docs/epexspot_prices.xlsx — add the "15-dec" spot-price column (inserted
# right before the "01-oct." block on the "Prix Spot" sheet) and append the
# 2025-12-13 / 2025-12-14 rows to the "Gaz" and "CO2" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Prix Spot": insert one column before EL (the start of the "01-oct."
# block) so everything from EL onward shifts right by one column
# (EL->EM, ..., FP->FQ). The new EL column gets the "15-dec" header and a
# "-" placeholder for every hourly data row.
# ---------------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")
$wsSpot.Columns("EL:EL").Insert()

$wsSpot.Range("EL1").Value = "15-dec"
for ($r = 2; $r -le 25; $r++) {
    $wsSpot.Range("EL$r").Value = "-"
}

# ---------------------------------------------------------------------------
# "Gaz": append the two new daily rows after the existing last row (169).
# Column A holds the date as plain text (matching the existing rows), so
# force a text format before assigning, then drop back to the sheet's
# normal (unstyled) cell style.
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$gazDates = $wsGaz.Range("A170:A171")
$gazDates.NumberFormat = "@"
$wsGaz.Range("A170").Value = "2025-12-13"
$wsGaz.Range("A171").Value = "2025-12-14"
$gazDates.Style = "Normal"
$wsGaz.Range("B170").Value = 26.075
$wsGaz.Range("B171").Value = 26.075

# ---------------------------------------------------------------------------
# "CO2": append the two new daily rows after the existing last row (170).
# ---------------------------------------------------------------------------
$wsCO2 = $wb.Worksheets.Item("CO2")
$co2Dates = $wsCO2.Range("A171:A172")
$co2Dates.NumberFormat = "@"
$wsCO2.Range("A171").Value = "2025-12-13"
$wsCO2.Range("A172").Value = "2025-12-14"
$co2Dates.Style = "Normal"
$wsCO2.Range("B171").Value = 84.09999999999999
$wsCO2.Range("B172").Value = 84.09999999999999
